# Corrected columns in scribe notes:
# Column B ("Handler") and column C ("Dog") had their underlying values
# reversed on every data row. On the "MStd" / "PremStd" sheets the header
# row already read correctly (B1=Handler, C1=Dog) so only the data rows
# (2..last) need to be swapped. On "MJWW" / "PremJWW" the header row
# itself was also reversed (B1=Dog, C1=Handler), so the swap there must
# include row 1 as well.

$wb = $excel.ActiveWorkbook

# Sheet name -> first row that needs the B/C swap applied.
$sheetsToFix = @{
    "MStd"    = 2
    "PremStd" = 2
    "MJWW"    = 1
    "PremJWW" = 1
}

foreach ($sheetName in $sheetsToFix.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
    $startRow = $sheetsToFix[$sheetName]

    if ($lastRow -ge $startRow) {
        $rowCount = $lastRow - $startRow + 1
        $bcRange = $ws.Range($ws.Cells.Item($startRow, 2), $ws.Cells.Item($lastRow, 3))
        $values = $bcRange.Value2

        $swapped = New-Object 'object[,]' $rowCount,2
        for ($i = 1; $i -le $rowCount; $i++) {
            $swapped[$i - 1, 0] = $values[$i, 2]
            $swapped[$i - 1, 1] = $values[$i, 1]
        }

        $bcRange.Value = $swapped
    }
}
